$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "2020" column (M) of data, mirroring the existing layout/styles of
# the adjacent "2019" column (L).

# Row 3 - thin separator row under the header block (style only, no value)
$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial(-4122)

# Row 4 - year header
$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("M4").Value = 2020

# Rows 5-13 - data values
$ws.Range("L5").Copy()
$ws.Range("M5").PasteSpecial(-4122)
$ws.Range("M5").Value = 34.377950588852634

$ws.Range("L6").Copy()
$ws.Range("M6").PasteSpecial(-4122)
$ws.Range("M6").Value = 4.8358243107925931

$ws.Range("L7").Copy()
$ws.Range("M7").PasteSpecial(-4122)
$ws.Range("M7").Value = 5.9543034993102522

$ws.Range("L8").Copy()
$ws.Range("M8").PasteSpecial(-4122)
$ws.Range("M8").Value = 51.21106605430419

$ws.Range("L9").Copy()
$ws.Range("M9").PasteSpecial(-4122)
$ws.Range("M9").Value = 27.156801192263725

$ws.Range("L10").Copy()
$ws.Range("M10").PasteSpecial(-4122)
$ws.Range("M10").Value = 0.94331159862228353

$ws.Range("L11").Copy()
$ws.Range("M11").PasteSpecial(-4122)
$ws.Range("M11").Value = 7.8509592890793316

$ws.Range("L12").Copy()
$ws.Range("M12").PasteSpecial(-4122)
$ws.Range("M12").Value = 64.733302669743793

$ws.Range("L13").Copy()
$ws.Range("M13").PasteSpecial(-4122)
$ws.Range("M13").Value = 97.67954817102779

# Row 14 - footer total
$ws.Range("L14").Copy()
$ws.Range("M14").PasteSpecial(-4122)
$ws.Range("M14").Value = 46.725153243037099

$excel.CutCopyMode = 0

# Move the active cell/selection to where the author last left it.
$ws.Range("L19").Select()
